# Updated cryptos list on Sun Sep 17 15:28:32 UTC 2023 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the coin rows that moved since the previous snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.729.29";   E = "  +0.49%  " },
    @{ Row = 3;  D = "1.640.60";    E = "  +0.05%  " },
    @{ Row = 4;  D = $null;         E = "  +0.36%  " },
    @{ Row = 5;  D = "217.80";      E = "  +1.64%  " },
    @{ Row = 6;  D = "0.504";       E = "  +0.10%  " },
    @{ Row = 7;  D = $null;         E = "  +0.43%  " },
    @{ Row = 8;  D = $null;         E = "  +0.36%  " },
    @{ Row = 9;  D = "0.0626";      E = "  +0.23%  " },
    @{ Row = 10; D = "19.14";       E = "  +0.28%  " },
    @{ Row = 11; D = $null;         E = "  +0.17%  " },
    @{ Row = 12; D = $null;         E = "  +0.02%  " },
    @{ Row = 13; D = "1.630.66";    E = "  -0.35%  " },
    @{ Row = 14; D = "4.16";        E = "  -0.31%  " },
    @{ Row = 15; D = $null;         E = "  -0.17%  " },
    @{ Row = 16; D = "64.68";       E = "  -0.15%  " },
    @{ Row = 17; D = "26.718.27";   E = $null },
    @{ Row = 18; D = "0.0₃0733";    E = "  -1.05%  " },
    @{ Row = 19; D = "215.30";      E = "  +0.15%  " },
    @{ Row = 20; D = $null;         E = "  +0.42%  " },
    @{ Row = 21; D = $null;         E = "  +0.93%  " },
    @{ Row = 22; D = "2.36";        E = "  +7.13%  " },
    @{ Row = 23; D = "6.23";        E = "  -0.09%  " },
    @{ Row = 24; D = "9.28";        E = "  -1.64%  " },
    @{ Row = 25; D = "145.36";      E = $null },
    @{ Row = 26; D = $null;         E = "  +0.31%  " },
    @{ Row = 27; D = $null;         E = "  -0.79%  " },
    @{ Row = 28; D = $null;         E = "  +0.80%  " },
    @{ Row = 29; D = "15.64";       E = "  -0.18%  " },
    @{ Row = 30; D = $null;         E = "  -0.69%  " },
    @{ Row = 31; D = $null;         E = "  +1.61%  " },
    @{ Row = 32; D = "3.39";        E = "  +1.04%  " },
    @{ Row = 33; D = "3.01";        E = "  +0.60%  " },
    @{ Row = 34; D = "1.288.40";    E = "  +1.00%  " },
    @{ Row = 35; D = $null;         E = "  +0.24%  " },
    @{ Row = 36; D = $null;         E = "  +1.32%  " },
    @{ Row = 37; D = $null;         E = "  -0.03%  " },
    @{ Row = 38; D = "0.538";       E = "  +1.51%  " },
    @{ Row = 39; D = $null;         E = "  -0.87%  " },
    @{ Row = 40; D = $null;         E = "  +0.51%  " },
    @{ Row = 41; D = "0.805";       E = "  -0.45%  " },
    @{ Row = 42; D = $null;         E = "  -1.07%  " },
    @{ Row = 43; D = $null;         E = "  -2.23%  " },
    @{ Row = 44; D = "1.779.39";    E = "  +0.03%  " },
    @{ Row = 45; D = "61.01";       E = "  +2.86%  " },
    @{ Row = 46; D = "91.75";       E = "  +0.62%  " },
    @{ Row = 47; D = $null;         E = "  +0.43%  " },
    @{ Row = 48; D = "0.0522";      E = "  +1.47%  " },
    @{ Row = 49; D = "7.64";        E = "  -1.08%  " },
    @{ Row = 50; D = $null;         E = "  +0.46%  " },
    @{ Row = 51; D = $null;         E = "  -0.03%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $looksNumeric = $u.D -match '^[+-]?[0-9]+(\.[0-9]+)?$'
        if ($looksNumeric) {
            # These "Price" strings are display text (e.g. "217.80"), not
            # real numbers — force the cell to Text so Excel doesn't
            # reinterpret/round them, then restore the default style so no
            # stray formatting is left behind.
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
